$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 13:59"

# Re-sort/relabel the province rows (values swapped between rows)
$ws.Range("A7").Value = "Bizkaia/Vizcaya"
$ws.Range("A8").Value = "Valencia/Valencia"
$ws.Range("A17").Value = "Malaga"
$ws.Range("A18").Value = "Gran Canaria"
$ws.Range("A25").Value = "Salamanca"
$ws.Range("A26").Value = "Sevilla"
$ws.Range("A27").Value = "Valladolid"
$ws.Range("A28").Value = "Granada"
$ws.Range("A29").Value = "Murcia"
$ws.Range("A30").Value = "Albacete"
$ws.Range("A31").Value = "Leon"
$ws.Range("A33").Value = "Tenerife"
$ws.Range("A34").Value = "Segovia"
$ws.Range("A35").Value = "Jaen"
$ws.Range("A36").Value = "Castello/Castellon"
$ws.Range("A37").Value = "Guadalajara"
$ws.Range("A38").Value = "Badajoz"
$ws.Range("A39").Value = "Soria"
$ws.Range("A42").Value = "Ourense"
$ws.Range("A45").Value = "Palencia"
$ws.Range("A46").Value = "Mallorca"
$ws.Range("A47").Value = "Cuenca"
$ws.Range("A48").Value = "Teruel"
$ws.Range("A49").Value = "Huesca"
$ws.Range("A50").Value = "Almeria"

# Updated case counts (Casos totales / Casos activos / Recuperados / Muertes)
$ws.Range("B7").Value = 2463
$ws.Range("C7").Value = 1023
$ws.Range("D7").Value = 1960
$ws.Range("E7").Value = 103
$ws.Range("B8").Value = 2263
$ws.Range("C8").Value = 68
$ws.Range("D8").Value = 2091
$ws.Range("E8").Value = 104
$ws.Range("B10").Value = 1779
$ws.Range("C10").Value = 1023
$ws.Range("D10").Value = 1250
$ws.Range("E10").Value = 109
$ws.Range("B17").Value = 1053
$ws.Range("C17").Value = 80
$ws.Range("D17").Value = 917
$ws.Range("E17").Value = 56
$ws.Range("B18").Value = 1025
$ws.Range("C18").Value = 25
$ws.Range("D18").Value = 964
$ws.Range("E18").Value = 36
$ws.Range("B24").Value = 894
$ws.Range("C24").Value = 1023
$ws.Range("D24").Value = 662
$ws.Range("E24").Value = 29
$ws.Range("B25").Value = 882
$ws.Range("C25").Value = 131
$ws.Range("D25").Value = 667
$ws.Range("E25").Value = 84
$ws.Range("B26").Value = 830
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 791
$ws.Range("E26").Value = 26
$ws.Range("B27").Value = 807
$ws.Range("C27").Value = 114
$ws.Range("D27").Value = 648
$ws.Range("E27").Value = 45
$ws.Range("B28").Value = 806
$ws.Range("C28").Value = 11
$ws.Range("D28").Value = 746
$ws.Range("E28").Value = 49
$ws.Range("B29").Value = 802
$ws.Range("C29").Value = 12
$ws.Range("D29").Value = 773
$ws.Range("E29").Value = 17
$ws.Range("B30").Value = 780
$ws.Range("C30").Value = 153
$ws.Range("D30").Value = 667
$ws.Range("E30").Value = 83
$ws.Range("B31").Value = 726
$ws.Range("C31").Value = 99
$ws.Range("D31").Value = 560
$ws.Range("E31").Value = 67
$ws.Range("B32").Value = 673
$ws.Range("C32").Value = 117
$ws.Range("D32").Value = 514
$ws.Range("E32").Value = 42
$ws.Range("B33").Value = 539
$ws.Range("C33").Value = 15
$ws.Range("D33").Value = 519
$ws.Range("E33").Value = 36
$ws.Range("B34").Value = 503
$ws.Range("C34").Value = 111
$ws.Range("D34").Value = 341
$ws.Range("E34").Value = 51
$ws.Range("B35").Value = 465
$ws.Range("C35").Value = 15
$ws.Range("D35").Value = 434
$ws.Range("E35").Value = 16
$ws.Range("B36").Value = 449
$ws.Range("C36").Value = 4
$ws.Range("D36").Value = 424
$ws.Range("E36").Value = 21
$ws.Range("B37").Value = 440
$ws.Range("C37").Value = 153
$ws.Range("D37").Value = 362
$ws.Range("E37").Value = 75
$ws.Range("B38").Value = 437
$ws.Range("C38").Value = 11
$ws.Range("D38").Value = 416
$ws.Range("E38").Value = 10
$ws.Range("B39").Value = 431
$ws.Range("C39").Value = 49
$ws.Range("D39").Value = 356
$ws.Range("E39").Value = 26
$ws.Range("B40").Value = 424
$ws.Range("C40").Value = 4
$ws.Range("D40").Value = 411
$ws.Range("E40").Value = 9
$ws.Range("B41").Value = 406
$ws.Range("C41").Value = 8
$ws.Range("D41").Value = 391
$ws.Range("E41").Value = 7
$ws.Range("B42").Value = 396
$ws.Range("C42").Value = 95
$ws.Range("D42").Value = 371
$ws.Range("E42").Value = 6
$ws.Range("B43").Value = 381
$ws.Range("C43").Value = 78
$ws.Range("D43").Value = 265
$ws.Range("E43").Value = 38
$ws.Range("B45").Value = 220
$ws.Range("C45").Value = 26
$ws.Range("D45").Value = 183
$ws.Range("E45").Value = 11
$ws.Range("B46").Value = 210
$ws.Range("C46").Value = 18
$ws.Range("D46").Value = 194
$ws.Range("E46").Value = 12
$ws.Range("B47").Value = 180
$ws.Range("C47").Value = 153
$ws.Range("D47").Value = 130
$ws.Range("E47").Value = 40
$ws.Range("B48").Value = 179
$ws.Range("C48").Value = 10
$ws.Range("D48").Value = 159
$ws.Range("E48").Value = 10
$ws.Range("B49").Value = 174
$ws.Range("C49").Value = 14
$ws.Range("D49").Value = 151
$ws.Range("E49").Value = 9
$ws.Range("B50").Value = 173
$ws.Range("C50").Value = 6
$ws.Range("D50").Value = 157
$ws.Range("E50").Value = 10
$ws.Range("B51").Value = 168
$ws.Range("C51").Value = 27
$ws.Range("D51").Value = 125
$ws.Range("E51").Value = 16
$ws.Range("B52").Value = 120
$ws.Range("D52").Value = 116
